# Generate Report for Handback
# Refresh the handoff/handback timestamps for file "12040c64-1c18-422e-ab12-e661cbc401b8"
# on the Overview sheet and on each locale sheet (zh-cn, de-de). The second file
# ("5d3e1a15-ee8c-4620-8cdf-9ca9e09c2bb7") is left untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-16 12:47:21"

# --- zh-cn sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-16 12:47:15"
$wsZhCn.Range("K2").Value = "2016-08-16 12:47:31"

# --- de-de sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-16 12:47:21"
$wsDeDe.Range("K2").Value = "2016-08-16 12:47:39"
